$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new rows are inserted right after the header, shifting all existing
# data rows down by 2 (old row 2 -> row 4, ... old row 21 -> row 23).
$ws.Range("A2:A3").EntireRow.Insert()
$ws.Range("A2:H3").ClearFormats()

$newTopRows = @(
    @(0, "walkingToRunning", 3.577567869787046, -6.830358225996321, 7.436059972857905, -0.3733085989952087, -2.073603868484497, 2.090381860733032),
    @(100, "walkingToRunning", 1.296618202114653, -9.790127622488438, 8.772651798817325, 0.200868934392929, -1.047221541404724, 2.184124946594238)
)

for ($i = 0; $i -lt $newTopRows.Length; $i++) {
    $row = $newTopRows[$i]
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# Eight brand-new rows are appended after the (now shifted) last data row (23).
$newBottomRows = @(
    @(2200, "walkingToRunning", -5.125987016035118, -3.317670953866559, -10.087697013307, -9.171860694885254, -10.49505233764648, -3.089466571807861),
    @(2300, "walkingToRunning", 12.1697812054039, -11.89667802884434, -0.7975602624165994, 5.832921981811523, -5.742907047271729, 5.470998287200928),
    @(2400, "walkingToRunning", 13.27590551955934, 1.014137889798754, 9.396237328566208, 0.7526758909225464, -7.854794025421143, 0.5755757093429565),
    @(2500, "walkingToRunning", -4.693403524588408, -26.48646446354427, 34.13960077485957, 8.307531356811523, 6.279134750366211, -0.9227187633514404),
    @(2600, "walkingToRunning", -13.05436339826206, -24.43262726167378, 17.33472581726497, -3.3331458568573, 3.669769525527954, -1.161337971687317),
    @(2700, "walkingToRunning", -15.37497096825698, 11.57301431993091, -14.26161232310758, -1.678790211677551, 2.915562152862549, 3.480551958084106),
    @(2800, "walkingToRunning", -9.968023679533061, -3.754608689092123, 4.53524044205469, -8.003265380859375, 1.636179566383362, -1.630586981773376),
    @(2900, "walkingToRunning", -3.008851450780502, -11.21099381420506, 19.69571330665882, 1.040297269821167, -2.371345281600952, 6.342917442321777)
)

for ($i = 0; $i -lt $newBottomRows.Length; $i++) {
    $row = $newBottomRows[$i]
    $r = 24 + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
